# Duplicate the text "Hello Git" to "Hello GitHello Git" and split the
# paragraph right after the new text (before the _GoBack bookmark) so the
# bookmark ends up alone in a new trailing paragraph.
#
# Using Find/Replace with a "^p" (paragraph mark) appended to the
# replacement text inserts the new paragraph break exactly at the point
# where the matched text ended - i.e. right before the bookmarkStart /
# bookmarkEnd that immediately followed "Hello Git" - which moves the
# bookmark into the newly created paragraph while leaving the original
# paragraph's own pPr/run untouched.
$d = $word.ActiveDocument
$d.Content.Find.Execute("Hello Git", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hello GitHello Git^p", 2)
